$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# hunk 0, row 112
$ws.Range("H112").Value = 2313.3
$ws.Range("J112").Value = 2408.7368
$ws.Range("L112").Value = 7226.2104
$ws.Range("N112").Value = -9442.2104
# hunk 1, row 132
$ws.Range("H132").Value = 127120.414
$ws.Range("I132").Value = 1871.463
$ws.Range("K132").Value = 5614.389
$ws.Range("M132").Value = -3084.389
# hunk 2, row 135
$ws.Range("H135").Value = 9091834
$ws.Range("I135").Value = 295.09375
$ws.Range("J135").Value = 21740932
$ws.Range("K135").Value = 2655.84375
$ws.Range("L135").Value = 195668388
$ws.Range("M135").Value = -120.84375
$ws.Range("N135").Value = -195673458
# hunk 3, row 137
$ws.Range("H137").Value = 22867.04
$ws.Range("I137").Value = 36538.855
$ws.Range("J137").Value = 5466.5454
$ws.Range("K137").Value = 109616.565
$ws.Range("L137").Value = 16399.6362
$ws.Range("M137").Value = -107066.565
$ws.Range("N137").Value = -21499.6362
# hunk 4, row 138
$ws.Range("H138").Value = 1620.49
$ws.Range("I138").Value = 823.8913
$ws.Range("J138").Value = 2299.074
$ws.Range("K138").Value = 2471.6739
$ws.Range("L138").Value = 6897.222
$ws.Range("M138").Value = 2668.3261
$ws.Range("N138").Value = -17177.222

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# hunk 5, row 61
$ws.Range("H61").Value = 942.12195
$ws.Range("I61").Value = 872.7714
$ws.Range("J61").Value = 1346.6666
$ws.Range("K61").Value = 872.7714
$ws.Range("L61").Value = 1346.6666
$ws.Range("M61").Value = -660.7714
$ws.Range("N61").Value = -1770.6666
# hunk 6, row 74
$ws.Range("H74").Value = 18332.225
$ws.Range("I74").Value = 24870.5
$ws.Range("K74").Value = 24870.5
$ws.Range("M74").Value = -23996.5
# hunk 7, row 77
$ws.Range("H77").Value = 18332.225
$ws.Range("I77").Value = 24870.5
$ws.Range("K77").Value = 124352.5
$ws.Range("M77").Value = -119984.5
# hunk 8, row 102
$ws.Range("H102").Value = 1300
$ws.Range("I102").Value = 1250
$ws.Range("K102").Value = 1250
$ws.Range("M102").Value = 372
# hunk 9, row 132
$ws.Range("H132").Value = 1651053.6
$ws.Range("I132").Value = 1891697.4
$ws.Range("J132").Value = 722856
$ws.Range("K132").Value = 5675092.199999999
$ws.Range("L132").Value = 2168568
$ws.Range("M132").Value = -5672562.199999999
$ws.Range("N132").Value = -2173628
# hunk 10, row 136
$ws.Range("H136").Value = 942.12195
$ws.Range("I136").Value = 872.7714
$ws.Range("J136").Value = 1346.6666
$ws.Range("K136").Value = 2618.3142
$ws.Range("L136").Value = 4039.9998
$ws.Range("M136").Value = -68.3141999999998
$ws.Range("N136").Value = -9139.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# hunk 11, row 29
$ws.Range("H29").Value = 27006
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 500
$ws.Range("M29").Value = -211
# hunk 12, row 107
$ws.Range("H107").Value = 788
$ws.Range("I107").Value = 788
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 788
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1132
$ws.Range("N107").ClearContents()
# hunk 13, row 134
$ws.Range("H134").Value = 32364.334
$ws.Range("I134").Value = 1496.6207
$ws.Range("J134").Value = 160244.86
$ws.Range("K134").Value = 4489.8621
$ws.Range("L134").Value = 480734.58
$ws.Range("M134").Value = -1954.8621
$ws.Range("N134").Value = -485804.58

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# hunk 14, row 16
$ws.Range("H16").Value = 798.0714
$ws.Range("I16").Value = 693.3333
$ws.Range("J16").Value = 986.6
$ws.Range("K16").Value = 693.3333
$ws.Range("L16").Value = 986.6
$ws.Range("M16").Value = -406.3333
$ws.Range("N16").Value = -1560.6
# hunk 15, row 58
$ws.Range("H58").Value = 897.58826
$ws.Range("I58").Value = 644.3684
$ws.Range("J58").Value = 1218.3334
$ws.Range("K58").Value = 644.3684
$ws.Range("L58").Value = 1218.3334
$ws.Range("M58").Value = -441.3684
$ws.Range("N58").Value = -1624.3334
# hunk 16, row 100
$ws.Range("H100").Value = 67000
$ws.Range("J100").Value = 67000
$ws.Range("L100").Value = 67000
$ws.Range("N100").Value = -69164
# hunk 17, row 113
$ws.Range("H113").Value = 798.0714
$ws.Range("I113").Value = 693.3333
$ws.Range("J113").Value = 986.6
$ws.Range("K113").Value = 693.3333
$ws.Range("L113").Value = 986.6
$ws.Range("M113").Value = 1476.6667
$ws.Range("N113").Value = -5326.6
# hunk 18, row 132
$ws.Range("H132").Value = 975.07275
$ws.Range("I132").Value = 624
$ws.Range("J132").Value = 1830.8125
$ws.Range("K132").Value = 1872
$ws.Range("L132").Value = 5492.4375
$ws.Range("M132").Value = 658
$ws.Range("N132").Value = -10552.4375
# hunk 19, row 134
$ws.Range("H134").Value = 948.3788
$ws.Range("I134").Value = 905.55316
$ws.Range("J134").Value = 1054.3158
$ws.Range("K134").Value = 2716.65948
$ws.Range("L134").Value = 3162.9474
$ws.Range("M134").Value = -181.6594800000003
$ws.Range("N134").Value = -8232.947400000001
# hunk 20, row 136
$ws.Range("H136").Value = 897.58826
$ws.Range("I136").Value = 644.3684
$ws.Range("J136").Value = 1218.3334
$ws.Range("K136").Value = 1933.1052
$ws.Range("L136").Value = 3655.0002
$ws.Range("M136").Value = 616.8948
$ws.Range("N136").Value = -8755.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# hunk 21, row 131
$ws.Range("H131").Value = 52084120
$ws.Range("I131").Value = 504
$ws.Range("J131").Value = 65790332
$ws.Range("K131").Value = 1512
$ws.Range("L131").Value = 197370996
$ws.Range("M131").Value = 3528
$ws.Range("N131").Value = -197381076

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# hunk 22, row 70
$ws.Range("H70").Value = 11116433
$ws.Range("I70").Value = 14290986
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 14290986
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -14290716
$ws.Range("N70").Value = -6040
# hunk 23, row 73
$ws.Range("H73").Value = 11116433
$ws.Range("I73").Value = 14290986
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 14290986
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -14290050
$ws.Range("N73").Value = -7372

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# hunk 24, row 136
$ws.Range("H136").Value = 279462.7
$ws.Range("I136").Value = 401265.1
$ws.Range("J136").Value = 2639.0908
$ws.Range("K136").Value = 1203795.3
$ws.Range("L136").Value = 7917.2724
$ws.Range("M136").Value = -1201245.3
$ws.Range("N136").Value = -13017.2724

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# hunk 25, row 93
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -24992
# hunk 26, row 123
$ws.Range("H123").Value = 36633.332
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 36633.332
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 36633.332
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -46433.332
# hunk 27, row 125
$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -34840
# hunk 28, row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# hunk 29, row 132
$ws.Range("H132").Value = 3822.842
$ws.Range("I132").Value = 997.2222
$ws.Range("J132").Value = 6365.9
$ws.Range("K132").Value = 2991.6666
$ws.Range("L132").Value = 19097.7
$ws.Range("M132").Value = -461.6666
$ws.Range("N132").Value = -24157.7
# hunk 30, row 136
$ws.Range("H136").Value = 1101317.2
$ws.Range("I136").Value = 1232444.9
$ws.Range("J136").Value = 625979.9399999999
$ws.Range("K136").Value = 3697334.7
$ws.Range("L136").Value = 1877939.82
$ws.Range("M136").Value = -3694784.7
$ws.Range("N136").Value = -1883039.82
